# Log workbook update: shift the "Files_Worked" (column C) entries down by one
# row (C7:C16 -> C8:C17), clearing C7 back to its default "empty day" look,
# and filling the newly used C17 cell. This reflects that the actual work
# originally logged against 2025-XX day 7 happened a day later, pushing every
# later day's note down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current Files_Worked text for rows 7 through 16 before we
# start overwriting any of them.
$vals = @()
for ($r = 7; $r -le 16; $r++) {
    $vals += $ws.Range("C$r").Value2
}

# Shift every captured value down by one row: old C16 -> C17, ..., old C7 -> C8.
# Walk from the bottom up so we never clobber a value before it's copied.
for ($i = $vals.Count - 1; $i -ge 0; $i--) {
    $destRow = 8 + $i
    $ws.Range("C$destRow").Value2 = $vals[$i]
}

# Row 7 no longer holds an entry, so clear it out and restore the plain
# "unfilled" cell styling (matching the look used higher up the sheet, e.g. C2).
$ws.Range("C7").ClearContents()
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the user's current working position: scrolled down slightly with
# C7 (the newly-cleared cell) selected.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
